$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "TCS Food" module: insert two new key/value rows right after the
# header block (before the existing "intro_1" row), pushing every
# following row down by two.
$ws.Rows("8:9").Insert()

$ws.Range("A8").Value = "tcs_foods"
$ws.Range("B8").Value = "TCS Foods\n(Time/Temperature Control for Safety)"

$ws.Range("A9").Value = "non_tcs_foods"
$ws.Range("B9").Value = "Non-TCS Foods"
